# Update "paises" (COVID countries) worksheet:
#  - refresh the "Datos actualizados" timestamp
#  - re-rank a handful of countries (their row labels swap) and refresh
#    their statistics to the newer snapshot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp cell (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 14:29"

# --- Row 4: Estados Unidos (name/rank unchanged, only numbers refresh) --
$ws.Range("B4").Value = 104277
$ws.Range("C4").Value = 151
$ws.Range("E4").Value = 100048

# --- Row 12: Suiza (name/rank unchanged, only numbers refresh) ----------
$ws.Range("B12").Value = 13377
$ws.Range("C12").Value = 449
$ws.Range("E12").Value = 11605
$ws.Range("G12").Value = 11
$ws.Range("H12").Value = 242

# --- Rows 13-15: Corea del Sur / Belgica / Paises Bajos swap places -----
# Row 13 becomes "Paises Bajos"
$ws.Range("A13").Value = "Paises Bajos"
$ws.Range("B13").Value = 9762
$ws.Range("C13").Value = 1159
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 9120
$ws.Range("F13").Value = 761
$ws.Range("G13").Value = 93
$ws.Range("H13").Value = 639

# Row 14 becomes "Corea del Sur"
$ws.Range("A14").Value = "Corea del Sur"
$ws.Range("B14").Value = 9478
$ws.Range("C14").Value = 146
$ws.Range("D14").Value = 4811
$ws.Range("E14").Value = 4523
$ws.Range("F14").Value = 59
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 144

# Row 15 becomes "Belgica"
$ws.Range("A15").Value = "Belgica"
$ws.Range("B15").Value = 9134
$ws.Range("C15").Value = 1850
$ws.Range("D15").Value = 1063
$ws.Range("E15").Value = 7718
$ws.Range("F15").Value = 789
$ws.Range("G15").Value = 64
$ws.Range("H15").Value = 353

# --- Rows 18-19: Canada / Portugal swap places ---------------------------
# Row 18 becomes "Portugal"
$ws.Range("A18").Value = "Portugal"
$ws.Range("B18").Value = 5170
$ws.Range("C18").Value = 902
$ws.Range("D18").Value = 43
$ws.Range("E18").Value = 5027
$ws.Range("F18").Value = 89
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 100

# Row 19 becomes "Canada"
$ws.Range("A19").Value = "Canada"
$ws.Range("B19").Value = 4757
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 354
$ws.Range("E19").Value = 4348
$ws.Range("F19").Value = 120
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 55

# --- Row 27: numbers refresh only ----------------------------------------
$ws.Range("E27").Value = 2135
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 65

# --- Rows 44-45: India / Islandia swap places -----------------------------
# Row 44 becomes "Islandia"
$ws.Range("A44").Value = "Islandia"
$ws.Range("B44").Value = 963
$ws.Range("C44").Value = 73
$ws.Range("D44").Value = 97
$ws.Range("E44").Value = 864
$ws.Range("F44").Value = 18
$ws.Range("H44").Value = 2

# Row 45 becomes "India"
$ws.Range("A45").Value = "India"
$ws.Range("B45").Value = 933
$ws.Range("C45").Value = 46
$ws.Range("D45").Value = 84
$ws.Range("E45").Value = 829
$ws.Range("F45").Value = 0
$ws.Range("H45").Value = 20
